$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Fix the title typo "CD/CD Workflow" -> "CI/CD Workflow" by replacing just
# the second character ("D" -> "I"), which matches how PowerPoint splits
# the original single run into three runs around the edited character.
$c = $tr.Characters(2, 1)
$c.Text = "I"

# Re-apply the (auto-fit driven) size PowerPoint computed for the shape
# after the text change.
$shp.Width = 849.10002
$shp.Height = 39.75
